$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.450.90'
$ws.Range("E2").Value = '  -1.32%  '

$ws.Range("D3").Value = '1.614.93'

$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.00'
$ws.Range("E5").Value = '  -1.24%  '

$ws.Range("E6").Value = '  -1.55%  '

$ws.Range("E7").Value = '  +0.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.75'
$ws.Range("E8").Value = '  -1.90%  '

$ws.Range("E9").Value = '  +0.22%  '

$ws.Range("E10").Value = '  -0.82%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0885'
$ws.Range("E11").Value = '  -0.67%  '

$ws.Range("D12").Value = '1.844.30'
$ws.Range("E12").Value = '  -2.12%  '

$ws.Range("D13").Value = '1.613.03'
$ws.Range("E13").Value = '  -2.61%  '

$ws.Range("E14").Value = '  -0.61%  '

$ws.Range("E15").Value = '  -3.16%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.03'
$ws.Range("E16").Value = '  +0.76%  '

$ws.Range("D17").Value = '27.428.23'
$ws.Range("E17").Value = '  -1.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.74'
$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("E19").Value = '  -1.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.50'
$ws.Range("E20").Value = '  -2.77%  '

$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("E22").Value = '  -0.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.15'
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.06'
$ws.Range("E24").Value = '  +5.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.96'
$ws.Range("E25").Value = '  -0.48%  '

$ws.Range("E26").Value = '  -1.90%  '

$ws.Range("E27").Value = '  -1.61%  '

$ws.Range("E28").Value = '  +0.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.51'
$ws.Range("E29").Value = '  -1.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.18'

$ws.Range("E31").Value = '  -1.43%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.26'
$ws.Range("E32").Value = '  -1.77%  '

$ws.Range("D33").Value = '1.470.99'
$ws.Range("E33").Value = '  +1.83%  '

$ws.Range("E34").Value = '  -3.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.53'
$ws.Range("E35").Value = '  -3.70%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.966'
$ws.Range("E36").Value = '  +9.91%  '

$ws.Range("E37").Value = '  -0.52%  '

$ws.Range("E38").Value = '  -0.86%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.555'
$ws.Range("E39").Value = '  -2.86%  '

$ws.Range("E40").Value = '  -3.06%  '

$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '66.80'
$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.984'
$ws.Range("E43").Value = '  -4.82%  '

$ws.Range("B44").Value = 'mCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.46'
$ws.Range("E44").Value = '  -0.21%  '

$ws.Range("E45").Value = '  -2.80%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.23'

$ws.Range("D47").Value = '1.755.37'
$ws.Range("E47").Value = '  -2.10%  '

$ws.Range("E48").Value = '  -0.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.68'
$ws.Range("E49").Value = '  +0.19%  '

$ws.Range("E50").Value = '  -2.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.100'
$ws.Range("E51").Value = '  +0.66%  '
